$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.492.27"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "'3.514.03"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'591.87"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").Value = "'134.81"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.488"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'7.52"
$ws.Range("E9").Value = "  +5.95%  "
$ws.Range("D10").Value = "'0.125"
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("D11").Value = "'0.386"
$ws.Range("E11").Value = "  +2.73%  "
$ws.Range("D12").Value = "'4.109.95"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("E13").Value = "  +1.73%  "
$ws.Range("D14").Value = "'0.0000182"
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").Value = "'3.508.55"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "'25.84"
$ws.Range("E16").Value = "  -4.84%  "
$ws.Range("D17").Value = "'64.452.05"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "'9.93"
$ws.Range("E18").Value = "  +1.86%  "
$ws.Range("D19").Value = "'5.77"
$ws.Range("E19").Value = "  +3.28%  "
$ws.Range("D20").Value = "'13.64"
$ws.Range("E20").Value = "  -1.65%  "
$ws.Range("D21").Value = "'394.32"
$ws.Range("E21").Value = "  +2.66%  "
$ws.Range("D22").Value = "'0.576"
$ws.Range("E22").Value = "  +1.52%  "
$ws.Range("D23").Value = "'3.652.47"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").Value = "'74.70"
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'0.0000118"
$ws.Range("E27").Value = "  +3.27%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").Value = "'7.42"
$ws.Range("E29").Value = "  -2.28%  "
$ws.Range("D30").Value = "'2.27"
$ws.Range("E30").Value = "  +2.32%  "
$ws.Range("D31").Value = "'8.30"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  -6.49%  "
$ws.Range("D33").Value = "'0.159"
$ws.Range("E33").Value = "  +8.66%  "
$ws.Range("D34").Value = "'3.537.36"
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'23.48"
$ws.Range("E36").Value = "  -0.42%  "
$ws.Range("D37").Value = "'5.36"
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("D38").Value = "'6.97"
$ws.Range("E38").Value = "  +1.44%  "
$ws.Range("D39").Value = "'1.56"
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("D40").Value = "'167.54"
$ws.Range("E40").Value = "  +2.01%  "
$ws.Range("D41").Value = "'0.0791"
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("D42").Value = "'0.811"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "'4.46"
$ws.Range("E44").Value = "  +1.55%  "
$ws.Range("D45").Value = "'25.06"
$ws.Range("E45").Value = "  -3.90%  "
$ws.Range("D46").Value = "'1.68"
$ws.Range("E46").Value = "  +3.22%  "
$ws.Range("E47").Value = "  -3.42%  "
$ws.Range("D48").Value = "'6.81"
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("D49").Value = "'2.386.02"
$ws.Range("E49").Value = "  -3.78%  "
$ws.Range("D50").Value = "'0.901"
$ws.Range("E50").Value = "  -1.81%  "
$ws.Range("D51").Value = "'0.220"
$ws.Range("E51").Value = "  +1.61%  "
